$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A2").Value = "95SY88G93C56"
$ws.Range("A3").Value = "HKDJB5BA6J6M"
$ws.Range("A4").Value = "04F5PQ59MWV6"

$ws.Range("A16").ClearContents()

$ws.Range("A10").Value = "C9DEXVFAR31A"
$ws.Range("A11").Value = "05ANDJ337D9B"
$ws.Range("A12").Value = "SAS5DZQK4GHR"
$ws.Range("A13").Value = "YC7CEVJY9735"
$ws.Range("A14").Value = "X36KP2Z510RZ"
$ws.Range("A15").Value = "MK83F9RSV97N"

$ws.Range("A5").Select()
